$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $searchText, $replaceText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed in paragraph" $paraIndex "for" $searchText
    }
    return $ok
}

# --- Paragraph 7: "PERIODO CONTRATADO: DE 30/10/2015 a 02/11/2015 (3 diárias)" ---
Replace-InParagraph 7 "30/10/2015" "12/11/2015"
Replace-InParagraph 7 "02/11/2015" "15/11/2015"

# --- Paragraph 14: dates "30 de Outubro de 2015" / "2 de Novembro de 2015" ---
# Replace the "2 de Novembro" one FIRST, since "12 de Novembro de 2015" (the
# replacement for the other date) contains "2 de Novembro de 2015" as a substring.
Replace-InParagraph 14 "2 de Novembro de 2015" "15 de Novembro de 2015"
Replace-InParagraph 14 "30 de Outubro de 2015" "12 de Novembro de 2015"

# --- Paragraph 17: rent value + amount in words + installment count ---
Replace-InParagraph 17 "R$ 2.850,00" "R$ 2.050,00"
Replace-InParagraph 17 "(DOIS MIL E OITOCENTOS E CINQ�ENTA REAIS" "(DOIS MIL E CINQ�ENTA REAIS"
Replace-InParagraph 17 "2 (DUAS)" "1 (UMA)"

# --- Paragraph 19: first installment (1a PARCELA) ---
Replace-InParagraph 19 " R$ 1.475,00" " R$ 2.150,00"
Replace-InParagraph 19 "UM MIL E QUATROCENTOS E SETENTA E CINCO REAIS" "DOIS MIL E CENTO E CINQ�ENTA REAIS"
Replace-InParagraph 19 "17 de Outubro de 2015" "10 de Novembro de 2015"

# --- Paragraph 46 / 92: signature dates ---
Replace-InParagraph 46 "17 de Outubro de 2015" "10 de Novembro de 2015"
Replace-InParagraph 92 "17 de Outubro de 2015" "10 de Novembro de 2015"

# --- Paragraphs 64 / 72: authorization period dates ---
Replace-InParagraph 64 "30/10/2015" "12/11/2015"
Replace-InParagraph 64 "02/11/2015" "15/11/2015"
Replace-InParagraph 72 "30/10/2015" "12/11/2015"
Replace-InParagraph 72 "02/11/2015" "15/11/2015"

# --- Paragraphs 89 / 90: entry / exit times ---
Replace-InParagraph 89 "30/10/2015" "12/11/2015"
Replace-InParagraph 90 "02/11/2015" "15/11/2015"

Write-Host "Done with text replacements"
